$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the column-level date-format style from column G (also clears
# the style left on the G1/G2 cells that inherited it).
$ws.Columns("G").EntireColumn.ClearFormats()

# Row 2 keeps only the Principle value (C2); everything else in that
# row -- the date, the booleans, the buy/sell prices -- is wiped out
# completely (value and formatting), including the leftover
# number-format style that used to sit on A2.
$ws.Range("A2").Clear()
$ws.Range("B2").Clear()
$ws.Range("D2").Clear()
$ws.Range("E2").Clear()
$ws.Range("F2").Clear()
$ws.Range("G2").Clear()

# Row 3 (the stray trailing "0") goes away entirely.
$ws.Range("A3:G3").Clear()

$wb.Save()
